$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 13677555
$ws.Range("I137").Value = 3473069
$ws.Range("J137").Value = 39917664
$ws.Range("K137").Value = 10419207
$ws.Range("L137").Value = 119752992
$ws.Range("M137").Value = -10416657
$ws.Range("N137").Value = -119758092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17000
$ws.Range("J24").Value = 17000
$ws.Range("L24").Value = 17000
$ws.Range("N24").Value = -17748
$ws.Range("H32").Value = 6973.19
$ws.Range("I32").Value = 6057.1934
$ws.Range("J32").Value = 19142.857
$ws.Range("K32").Value = 6057.1934
$ws.Range("L32").Value = 19142.857
$ws.Range("M32").Value = -5770.1934
$ws.Range("N32").Value = -19716.857
$ws.Range("H37").Value = 9209.632
$ws.Range("J37").Value = 14045.2
$ws.Range("L37").Value = 14045.2
$ws.Range("N37").Value = -14591.2
$ws.Range("H61").Value = 3243071
$ws.Range("I61").Value = 1667897.2
$ws.Range("J61").Value = 9806295
$ws.Range("K61").Value = 1667897.2
$ws.Range("L61").Value = 9806295
$ws.Range("M61").Value = -1667685.2
$ws.Range("N61").Value = -9806719
$ws.Range("H62").Value = 8000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 8000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").Value = 24000
$ws.Range("N65").Value = -30240
$ws.Range("H100").Value = 17000
$ws.Range("J100").Value = 17000
$ws.Range("L100").Value = 17000
$ws.Range("N100").Value = -19164
$ws.Range("H136").Value = 3243071
$ws.Range("I136").Value = 1667897.2
$ws.Range("J136").Value = 9806295
$ws.Range("K136").Value = 5003691.6
$ws.Range("L136").Value = 29418885
$ws.Range("M136").Value = -5001141.6
$ws.Range("N136").Value = -29423985

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H86").Value = 1994
$ws.Range("I86").Value = 1994
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1994
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -871
$ws.Range("H89").Value = 1994
$ws.Range("I89").Value = 1994
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9970
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -4354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2528663.2
$ws.Range("I31").Value = 3473290.5
$ws.Range("K31").Value = 3473290.5
$ws.Range("M31").Value = -3472995.5
$ws.Range("H34").Value = 2528663.2
$ws.Range("I34").Value = 3473290.5
$ws.Range("K34").Value = 3473290.5
$ws.Range("M34").Value = -3473088.5
$ws.Range("H50").Value = 12624.111
$ws.Range("J50").Value = 12624.111
$ws.Range("L50").Value = 12624.111
$ws.Range("N50").Value = -13874.111
$ws.Range("H51").Value = 28330.1
$ws.Range("J51").Value = 10412.625
$ws.Range("L51").Value = 10412.625
$ws.Range("N51").Value = -11884.625
$ws.Range("H58").Value = 2532808.2
$ws.Range("I58").Value = 12276.223
$ws.Range("J58").Value = 5053340
$ws.Range("K58").Value = 12276.223
$ws.Range("L58").Value = 5053340
$ws.Range("M58").Value = -12073.223
$ws.Range("N58").Value = -5053746
$ws.Range("H59").Value = 17997.143
$ws.Range("J59").Value = 17997.143
$ws.Range("L59").Value = 17997.143
$ws.Range("N59").Value = -20287.143
$ws.Range("H60").Value = 17692.75
$ws.Range("J60").Value = 10330.429
$ws.Range("L60").Value = 10330.429
$ws.Range("N60").Value = -11352.429
$ws.Range("H61").Value = 28330.1
$ws.Range("J61").Value = 10412.625
$ws.Range("L61").Value = 10412.625
$ws.Range("N61").Value = -11108.625
$ws.Range("H68").Value = 17255
$ws.Range("J68").Value = 17255
$ws.Range("L68").Value = 17255
$ws.Range("N68").Value = -18753
$ws.Range("H71").Value = 17255
$ws.Range("J71").Value = 17255
$ws.Range("L71").Value = 51765
$ws.Range("N71").Value = -59253
$ws.Range("H74").Value = 17233.928
$ws.Range("J74").Value = 18460.77
$ws.Range("L74").Value = 18460.77
$ws.Range("N74").Value = -20208.77
$ws.Range("H77").Value = 17233.928
$ws.Range("J77").Value = 18460.77
$ws.Range("L77").Value = 55382.31
$ws.Range("N77").Value = -64118.31
$ws.Range("H92").Value = 33612.5
$ws.Range("J92").Value = 33612.5
$ws.Range("L92").Value = 33612.5
$ws.Range("N92").Value = -38604.5
$ws.Range("H94").Value = 35721484
$ws.Range("I94").Value = 1703
$ws.Range("J94").Value = 50009396
$ws.Range("K94").Value = 1703
$ws.Range("L94").Value = 50009396
$ws.Range("M94").Value = -1252
$ws.Range("N94").Value = -50010298
$ws.Range("H96").Value = 16841.334
$ws.Range("J96").Value = 16841.334
$ws.Range("L96").Value = 16841.334
$ws.Range("N96").Value = -22333.334
$ws.Range("H134").Value = 1145409.4
$ws.Range("I134").Value = 1915.8636
$ws.Range("J134").Value = 3080552
$ws.Range("K134").Value = 5747.5908
$ws.Range("L134").Value = 9241656
$ws.Range("M134").Value = -3212.5908
$ws.Range("N134").Value = -9246726
$ws.Range("H136").Value = 2532808.2
$ws.Range("I136").Value = 12276.223
$ws.Range("J136").Value = 5053340
$ws.Range("K136").Value = 36828.669
$ws.Range("L136").Value = 15160020
$ws.Range("M136").Value = -34278.669
$ws.Range("N136").Value = -15165120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2139.7144
$ws.Range("I132").Value = 1373.5
$ws.Range("J132").Value = 2320
$ws.Range("K132").Value = 12361.5
$ws.Range("L132").Value = 20880
$ws.Range("M132").Value = -9831.5
$ws.Range("N132").Value = -25940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4616.6665
$ws.Range("I41").Value = 833.3333
$ws.Range("J41").Value = 8400
$ws.Range("K41").Value = 833.3333
$ws.Range("L41").Value = 8400
$ws.Range("M41").Value = -478.3333
$ws.Range("N41").Value = -9110
$ws.Range("H70").Value = 7338399.5
$ws.Range("I70").Value = 2845275
$ws.Range("J70").Value = 23813190
$ws.Range("K70").Value = 2845275
$ws.Range("L70").Value = 23813190
$ws.Range("M70").Value = -2845005
$ws.Range("N70").Value = -23813730
$ws.Range("H73").Value = 7338399.5
$ws.Range("I73").Value = 2845275
$ws.Range("J73").Value = 23813190
$ws.Range("K73").Value = 2845275
$ws.Range("L73").Value = 23813190
$ws.Range("M73").Value = -2844339
$ws.Range("N73").Value = -23815062
$ws.Range("H107").Value = 341.23077
$ws.Range("I107").Value = 114.333336
$ws.Range("K107").Value = 114.333336
$ws.Range("M107").Value = 1805.666664
$ws.Range("H113").Value = 20576
$ws.Range("I113").Value = 1040.4445
$ws.Range("J113").Value = 55740
$ws.Range("K113").Value = 1040.4445
$ws.Range("L113").Value = 55740
$ws.Range("M113").Value = 1129.5555
$ws.Range("N113").Value = -60080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 15077
$ws.Range("I93").Value = 3415.9
$ws.Range("J93").Value = 28033.777
$ws.Range("K93").Value = 3415.9
$ws.Range("L93").Value = 28033.777
$ws.Range("M93").Value = -2167.9
$ws.Range("N93").Value = -30529.777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 21055
$ws.Range("J48").Value = 21055
$ws.Range("L48").Value = 21055
$ws.Range("N48").Value = -22193
$ws.Range("H69").Value = 11333.333
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16498
$ws.Range("H72").Value = 11333.333
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -52488
